$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKKPS134-001")
$ws2 = $wb.Worksheets.Item("DPLKKPS134-002")

# --- Sheet "DPLKKPS134-001": No. Register changes from M13220800000030 to M13220800000028
# (status stays "0 : Kembalikan ke Register")
$ws1.Range("N2").Value = "M13220800000028"
$ws1.Range("F2").Value = "Username : 31816;`nPassword : bni1234;`nRole : 09 - Penyelia Settlement;`nNo. Register : M13220800000028;`nStatus Verifikasi : 0 : Kembalikan ke Register;`nKeterangan Verifikasi : KEP.TRX.436 kembalikan ke Register"

# --- Sheet "DPLKKPS134-002": No. Register changes from M13220800000030 to M13220800000027
# (status stays "1 : Setuju")
$ws2.Range("N2").Value = "M13220800000027"
$ws2.Range("F2").Value = "Username : 31816;`nPassword : bni1234;`nRole : 09 - Penyelia Settlement;`nNo. Register : M13220800000027;`nStatus Verifikasi : 1 : Setuju;`nKeterangan Verifikasi : KEP.TRX.436 Setuju"

# --- View state changes ---
# Selection on "DPLKKPS134-002" moves from G2 to F2 (set first, since it is still active here).
$ws2.Range("F2").Select()

# Selection on "DPLKKPS134-001" moves from Q2 to G2, and it becomes the active (tab-selected)
# sheet -- select it last so it ends up as the workbook's active sheet.
$ws1.Range("G2").Select()
